$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of mileage data (row 8 and row 9)
$ws.Range("A8").Value = 43742
$ws.Range("B8").Value = 192

$ws.Range("A9").Value = 43743
$ws.Range("B9").Value = 240

# Update the selected cell to B10, matching the diff's selection change
$ws.Range("B10").Select()
